# Add a "Violations_Summary" sheet after Sheet1 summarizing rule violations
# with COUNTIF formulas against (currently empty) Sheet1 columns K:O.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Violations_Summary"

# Header row
$ws.Range("A1").Value = "A (Violation Type)"
$ws.Range("B1").Value = "`tB (Count)"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").HorizontalAlignment = -4108
$ws.Range("A1:B1").VerticalAlignment = -4108
$ws.Range("A1:B1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 43.5

# Row 2 - MFA
$ws.Range("A2").Value = "MFA"
$ws.Range("B2").Formula = "=COUNTIF(Sheet1!K2:K31, 1)"

# Row 3 - Password Age > 90 Days
$ws.Range("A3").Value = "Password Age > 90 Days"
$ws.Range("B3").Formula = "=COUNTIF(Sheet1!L2:L31, 1)"
$ws.Rows.Item(3).RowHeight = 43.5

# Row 4 - Inactive > 60 Days
$ws.Range("A4").Value = "Inactive > 60 Days"
$ws.Range("B4").Formula = "=COUNTIF(Sheet1!M2:M31, 1)"
$ws.Rows.Item(4).RowHeight = 29

# Row 5 - Terminated but Active
$ws.Range("A5").Value = "Terminated but Active"
$ws.Range("B5").Formula = "=COUNTIF(Sheet1!N2:N31, 1)"
$ws.Rows.Item(5).RowHeight = 43.5

# Row 6 - Admin Review Overdue
$ws.Range("A6").Value = "Admin Review Overdue"
$ws.Range("B6").Formula = "=COUNTIF(Sheet1!O2:O31, 1)"
$ws.Rows.Item(6).RowHeight = 43.5

# Body formatting: vertical-center + wrap for A2:A6, same + Arial Unicode MS 10pt for B2:B6
$ws.Range("A2:A6").VerticalAlignment = -4108
$ws.Range("A2:A6").WrapText = $true

$ws.Range("B2:B6").VerticalAlignment = -4108
$ws.Range("B2:B6").WrapText = $true
$ws.Range("B2:B6").Font.Name = "Arial Unicode MS"
$ws.Range("B2:B6").Font.Size = 10

$ws.Activate()
$ws.Range("E5").Select()
